$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @'
questions = [
    {
        "title": "You are a first-line manager at a high-end retail store. One of your sales associates is struggling to handle an upset customer who is dissatisfied with a recent purchase. The customer is demanding a full refund and is becoming increasingly agitated.How should you support your staff in this situation?",
        "ques_type": 2,
        "options": [
            "Step in and calmly address the customer's concerns.",
            "Step in and offer the customer store credit instead of a refund.",
            "Allow your associate to handle the situation on their own as a learning experience.",
            "Call for security to remove the upset customer from the store."
        ],
        "score": "Step in and calmly address the customer's concerns."
    },
    {
        "title": "You are a first-line manager at a high-end clothing retail store. Recently, several customers have expressed dissatisfaction with the limited variety of women's dresses. Your goal is to use this feedback to enhance customer satisfaction and boost dress sales.What should be your first course of action?",
        "ques_type": 2,
        "options": [
            "Buy a wide range of dresses in different styles and sizes.",
            "Analyze the customer feedback to pinpoint which styles, colors, or sizes of dresses are requested.",
            "Carry out a survey to gather additional feedback on specific styles of dresses customers would like to see.",
            "Implement a discount strategy on existing dresses to draw in customers."
        ],
        "score": "Analyze the customer feedback to pinpoint which styles, colors, or sizes of dresses are requested."
    },
    {
        "title": "You are a first-line manager at a high-end fashion retail store. You need to ensure a new collection is displayed in four locations across the store in a manner that will optimize sales.Which actions should you take?",
        "ques_type": 15,
        "options": [
            "Include as many items of the collection as possible in the four displays.",
            "Integrate the new collection with current collections in four existing displays.",
            "Delegate each of the four displays to a different sales associate to create.",
            "Display the collection in four high-traffic areas of the store.",
            "Create the four displays around a single unifying theme.",
            "Place all items in the four displays at the same height and depth."
        ],
        "score": [
            "Display the collection in four high-traffic areas of the store.",
            "Create the four displays around a single unifying theme."
        ]
    },
    {
        "title": "You are a first-line manager at a retail store, leading a team with diverse backgrounds and experience levels. Your goal is to foster strong relationships and trust among your staff.Which action should you take to achieve this?",
        "ques_type": 2,
        "options": [
            "Use email communication to address team members' concerns and questions.",
            "Interact with team members in formal meetings and performance evaluations only.",
            "Encourage team members to resolve their issues independently of management.",
            "Hold regular one-on-one meetings with each team member."
        ],
        "score": "Hold regular one-on-one meetings with each team member."
    }
]
'@

$ws.Range("A2").ClearContents()
$ws.Range("A1").Value = $text
$ws.Range("A1").Style = "Normal"
$ws.Rows(1).AutoFit()
